$d = $word.ActiveDocument

# 1. "専用サプライヤー" -> "専属サプライヤー" and add space before "ドリンク"
$d.Content.Find.Execute(
    "Northwind Traders は、Contoso のソフトドリンクとジュースの専用サプライヤーです。",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Northwind Traders は、Contoso のソフト ドリンクとジュースの専属サプライヤーです。", 2)

# 2. "一律料金" -> "遅延料金"
$d.Content.Find.Execute(
    "一律料金", $true, $false, $false, $false, $false, $true, 1, $false,
    "遅延料金", 2)

# 3. "1 か月あたり $22.5" -> "1 か月あたり 1.5%"
$d.Content.Find.Execute(
    "1 か月あたり `$22.5", $true, $false, $false, $false, $false, $true, 1, $false,
    "1 か月あたり 1.5%", 2)

# 4. "早期支払い割引" -> "早期支払割引"
$d.Content.Find.Execute(
    "早期支払い割引", $true, $false, $false, $false, $false, $true, 1, $false,
    "早期支払割引", 2)

# 5. "10 日以内の 2% 割引" -> "10 日以内 2% 割引"
$d.Content.Find.Execute(
    "10 日以内の 2% 割引", $true, $false, $false, $false, $false, $true, 1, $false,
    "10 日以内 2% 割引", 2)

# 6. "契約は別の年に自動的に更新されます。" -> "契約はもう 1 年自動的に更新されます。"
$d.Content.Find.Execute(
    "いずれかの当事者が有効期限の少なくとも 30 日前に終了の書面による通知を行わない限り、契約は別の年に自動的に更新されます。",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "いずれかの当事者が有効期限の少なくとも 30 日前に終了の書面による通知を行わない限り、契約はもう 1 年自動的に更新されます。", 2)

# 7. "最小注文金額" -> "最小注文量"
$d.Content.Find.Execute(
    "最小注文金額", $true, $false, $false, $false, $false, $true, 1, $false,
    "最小注文量", 2)

# 8. "1 か月あたり 100 リリース" -> "1 か月あたり 100 ケース"
$d.Content.Find.Execute(
    "1 か月あたり 100 リリース", $true, $false, $false, $false, $false, $true, 1, $false,
    "1 か月あたり 100 ケース", 2)

# 9. "最大注文金額" -> "最大注文量"
$d.Content.Find.Execute(
    "最大注文金額", $true, $false, $false, $false, $false, $true, 1, $false,
    "最大注文量", 2)

# 10. "20 時間/月" -> "1 か月あたり 500 ケース"
$d.Content.Find.Execute(
    "20 時間/月", $true, $false, $false, $false, $false, $true, 1, $false,
    "1 か月あたり 500 ケース", 2)
